# Apply data refresh to the Leve profit tables across all 8 job sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(113, 8).Value = 2651.7856  # H113
$ws.Cells.Item(113, 9).Value = 2385.4167  # I113
$ws.Cells.Item(113, 10).Value = 4250  # J113
$ws.Cells.Item(113, 11).Value = 2385.4167  # K113
$ws.Cells.Item(113, 12).Value = 4250  # L113
$ws.Cells.Item(113, 13).Value = 868.5832999999998  # M113
$ws.Cells.Item(113, 14).Value = -10758  # N113

$ws.Cells.Item(116, 8).Value = 21102  # H116
$ws.Cells.Item(116, 9).Value = 21102  # I116
$ws.Cells.Item(116, 10).Value = 0  # J116
$ws.Cells.Item(116, 11).Value = 21102  # K116
$ws.Cells.Item(116, 12).Value = 0  # L116
$ws.Cells.Item(116, 13).Value = -17660  # M116
$ws.Cells.Item(116, 14).ClearContents()  # N116

$ws.Cells.Item(132, 8).Value = 53566.42  # H132
$ws.Cells.Item(132, 9).Value = 56486.777  # I132
$ws.Cells.Item(132, 10).Value = 1000  # J132
$ws.Cells.Item(132, 11).Value = 169460.331  # K132
$ws.Cells.Item(132, 12).Value = 3000  # L132
$ws.Cells.Item(132, 13).Value = -166930.331  # M132
$ws.Cells.Item(132, 14).Value = -8060  # N132

$ws.Cells.Item(136, 8).Value = 0  # H136
$ws.Cells.Item(136, 10).Value = 0  # J136
$ws.Cells.Item(136, 12).Value = 0  # L136
$ws.Cells.Item(136, 14).ClearContents()  # N136

$ws.Cells.Item(137, 8).Value = 901.4286  # H137
$ws.Cells.Item(137, 9).Value = 776.6667  # I137
$ws.Cells.Item(137, 11).Value = 2330.0001  # K137
$ws.Cells.Item(137, 13).Value = 219.9998999999998  # M137

$ws.Cells.Item(138, 8).Value = 3145.0151  # H138
$ws.Cells.Item(138, 9).Value = 1993.7059  # I138
$ws.Cells.Item(138, 10).Value = 3544.449  # J138
$ws.Cells.Item(138, 11).Value = 5981.1177  # K138
$ws.Cells.Item(138, 12).Value = 10633.347  # L138
$ws.Cells.Item(138, 13).Value = -841.1176999999998  # M138
$ws.Cells.Item(138, 14).Value = -20913.347  # N138

$ws.Cells.Item(141, 8).Value = 5134.4116  # H141
$ws.Cells.Item(141, 9).Value = 5865.4165  # I141
$ws.Cells.Item(141, 10).Value = 3380  # J141
$ws.Cells.Item(141, 11).Value = 17596.2495  # K141
$ws.Cells.Item(141, 12).Value = 10140  # L141
$ws.Cells.Item(141, 13).Value = -12416.2495  # M141
$ws.Cells.Item(141, 14).Value = -20500  # N141

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(74, 8).Value = 976.125  # H74
$ws.Cells.Item(74, 9).Value = 957.25  # I74
$ws.Cells.Item(74, 11).Value = 957.25  # K74
$ws.Cells.Item(74, 13).Value = -83.25  # M74

$ws.Cells.Item(77, 8).Value = 976.125  # H77
$ws.Cells.Item(77, 9).Value = 957.25  # I77
$ws.Cells.Item(77, 11).Value = 4786.25  # K77
$ws.Cells.Item(77, 13).Value = -418.25  # M77

$ws.Cells.Item(122, 8).Value = 27678.104  # H122
$ws.Cells.Item(122, 9).Value = 2033.0714  # I122
$ws.Cells.Item(122, 10).Value = 92956.37  # J122
$ws.Cells.Item(122, 11).Value = 6099.2142  # K122
$ws.Cells.Item(122, 12).Value = 278869.11  # L122
$ws.Cells.Item(122, 13).Value = -3649.2142  # M122
$ws.Cells.Item(122, 14).Value = -283769.11  # N122

$ws.Cells.Item(132, 8).Value = 17876780  # H132
$ws.Cells.Item(132, 9).Value = 27779364  # I132
$ws.Cells.Item(132, 10).Value = 52128.9  # J132
$ws.Cells.Item(132, 11).Value = 83338092  # K132
$ws.Cells.Item(132, 12).Value = 156386.7  # L132
$ws.Cells.Item(132, 13).Value = -83335562  # M132
$ws.Cells.Item(132, 14).Value = -161446.7  # N132

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 587.5  # H80
$ws.Cells.Item(80, 10).Value = 423.76923  # J80
$ws.Cells.Item(80, 12).Value = 423.76923  # L80
$ws.Cells.Item(80, 14).Value = -2419.76923  # N80

$ws.Cells.Item(83, 8).Value = 587.5  # H83
$ws.Cells.Item(83, 10).Value = 423.76923  # J83
$ws.Cells.Item(83, 12).Value = 2118.84615  # L83
$ws.Cells.Item(83, 14).Value = -12102.84615  # N83

$ws.Cells.Item(134, 8).Value = 4083.83  # H134
$ws.Cells.Item(134, 9).Value = 1415.3658  # I134
$ws.Cells.Item(134, 11).Value = 4246.097400000001  # K134
$ws.Cells.Item(134, 13).Value = -1711.097400000001  # M134

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 10005173  # H31
$ws.Cells.Item(31, 9).Value = 17859842  # I31
$ws.Cells.Item(31, 10).Value = 8319.454  # J31
$ws.Cells.Item(31, 11).Value = 17859842  # K31
$ws.Cells.Item(31, 12).Value = 8319.454  # L31
$ws.Cells.Item(31, 13).Value = -17859547  # M31
$ws.Cells.Item(31, 14).Value = -8909.454  # N31

$ws.Cells.Item(34, 8).Value = 10005173  # H34
$ws.Cells.Item(34, 9).Value = 17859842  # I34
$ws.Cells.Item(34, 10).Value = 8319.454  # J34
$ws.Cells.Item(34, 11).Value = 17859842  # K34
$ws.Cells.Item(34, 12).Value = 8319.454  # L34
$ws.Cells.Item(34, 13).Value = -17859640  # M34
$ws.Cells.Item(34, 14).Value = -8723.454  # N34

$ws.Cells.Item(58, 8).Value = 938.2963  # H58
$ws.Cells.Item(58, 9).Value = 860.8  # I58
$ws.Cells.Item(58, 10).Value = 1035.1666  # J58
$ws.Cells.Item(58, 11).Value = 860.8  # K58
$ws.Cells.Item(58, 12).Value = 1035.1666  # L58
$ws.Cells.Item(58, 13).Value = -657.8  # M58
$ws.Cells.Item(58, 14).Value = -1441.1666  # N58

$ws.Cells.Item(132, 8).Value = 38128.465  # H132
$ws.Cells.Item(132, 9).Value = 1618.1364  # I132
$ws.Cells.Item(132, 10).Value = 171999.67  # J132
$ws.Cells.Item(132, 11).Value = 4854.4092  # K132
$ws.Cells.Item(132, 12).Value = 515999.01  # L132
$ws.Cells.Item(132, 13).Value = -2324.4092  # M132
$ws.Cells.Item(132, 14).Value = -521059.01  # N132

$ws.Cells.Item(134, 8).Value = 1779.9706  # H134
$ws.Cells.Item(134, 9).Value = 1301.8846  # I134
$ws.Cells.Item(134, 10).Value = 3333.75  # J134
$ws.Cells.Item(134, 11).Value = 3905.6538  # K134
$ws.Cells.Item(134, 12).Value = 10001.25  # L134
$ws.Cells.Item(134, 13).Value = -1370.6538  # M134
$ws.Cells.Item(134, 14).Value = -15071.25  # N134

$ws.Cells.Item(136, 8).Value = 938.2963  # H136
$ws.Cells.Item(136, 9).Value = 860.8  # I136
$ws.Cells.Item(136, 10).Value = 1035.1666  # J136
$ws.Cells.Item(136, 11).Value = 2582.4  # K136
$ws.Cells.Item(136, 12).Value = 3105.4998  # L136
$ws.Cells.Item(136, 13).Value = -32.39999999999964  # M136
$ws.Cells.Item(136, 14).Value = -8205.4998  # N136

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(98, 8).Value = 996.25  # H98
$ws.Cells.Item(98, 10).Value = 998.3333  # J98
$ws.Cells.Item(98, 12).Value = 2994.9999  # L98
$ws.Cells.Item(98, 14).Value = -5990.9999  # N98

$ws.Cells.Item(107, 8).Value = 250.13043  # H107
$ws.Cells.Item(107, 10).Value = 252.10811  # J107
$ws.Cells.Item(107, 12).Value = 756.32433  # L107
$ws.Cells.Item(107, 14).Value = -4596.32433  # N107

$ws.Cells.Item(113, 8).Value = 1021.5211  # H113
$ws.Cells.Item(113, 9).Value = 906.5  # I113
$ws.Cells.Item(113, 10).Value = 1024.8551  # J113
$ws.Cells.Item(113, 11).Value = 2719.5  # K113
$ws.Cells.Item(113, 12).Value = 3074.5653  # L113
$ws.Cells.Item(113, 13).Value = -549.5  # M113
$ws.Cells.Item(113, 14).Value = -7414.5653  # N113

$ws.Cells.Item(131, 8).Value = 6579955  # H131
$ws.Cells.Item(131, 9).Value = 1393.3334  # I131
$ws.Cells.Item(131, 10).Value = 8197634  # J131
$ws.Cells.Item(131, 11).Value = 4180.0002  # K131
$ws.Cells.Item(131, 12).Value = 24592902  # L131
$ws.Cells.Item(131, 13).Value = 859.9997999999996  # M131
$ws.Cells.Item(131, 14).Value = -24602982  # N131

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(102, 8).Value = 2173.1428  # H102
$ws.Cells.Item(102, 9).Value = 1200  # I102
$ws.Cells.Item(102, 10).Value = 2903  # J102
$ws.Cells.Item(102, 11).Value = 1200  # K102
$ws.Cells.Item(102, 12).Value = 2903  # L102
$ws.Cells.Item(102, 13).Value = 422  # M102
$ws.Cells.Item(102, 14).Value = -6147  # N102

$ws.Cells.Item(132, 8).Value = 208653.58  # H132
$ws.Cells.Item(132, 9).Value = 29342.285  # I132
$ws.Cells.Item(132, 10).Value = 557314.4399999999  # J132
$ws.Cells.Item(132, 11).Value = 88026.855  # K132
$ws.Cells.Item(132, 12).Value = 1671943.32  # L132
$ws.Cells.Item(132, 13).Value = -85496.855  # M132
$ws.Cells.Item(132, 14).Value = -1677003.32  # N132

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 2658.1  # H7
$ws.Cells.Item(7, 9).Value = 1449.8334  # I7
$ws.Cells.Item(7, 10).Value = 3175.9285  # J7
$ws.Cells.Item(7, 11).Value = 1449.8334  # K7
$ws.Cells.Item(7, 12).Value = 3175.9285  # L7
$ws.Cells.Item(7, 13).Value = -1337.8334  # M7
$ws.Cells.Item(7, 14).Value = -3399.9285  # N7

$ws.Cells.Item(94, 8).Value = 0  # H94
$ws.Cells.Item(94, 10).Value = 0  # J94
$ws.Cells.Item(94, 12).Value = 0  # L94
$ws.Cells.Item(94, 14).ClearContents()  # N94

$ws.Cells.Item(122, 8).Value = 2714.7856  # H122
$ws.Cells.Item(122, 9).Value = 2572.5715  # I122
$ws.Cells.Item(122, 10).Value = 2857  # J122
$ws.Cells.Item(122, 11).Value = 7717.7145  # K122
$ws.Cells.Item(122, 12).Value = 8571  # L122
$ws.Cells.Item(122, 13).Value = -5267.7145  # M122
$ws.Cells.Item(122, 14).Value = -13471  # N122

$ws.Cells.Item(126, 8).Value = 2658.1  # H126
$ws.Cells.Item(126, 9).Value = 1449.8334  # I126
$ws.Cells.Item(126, 10).Value = 3175.9285  # J126
$ws.Cells.Item(126, 11).Value = 4349.5002  # K126
$ws.Cells.Item(126, 12).Value = 9527.7855  # L126
$ws.Cells.Item(126, 13).Value = -1879.5002  # M126
$ws.Cells.Item(126, 14).Value = -14467.7855  # N126

$ws.Cells.Item(132, 8).Value = 23793.334  # H132
$ws.Cells.Item(132, 10).Value = 1587.9474  # J132
$ws.Cells.Item(132, 12).Value = 4763.8422  # L132
$ws.Cells.Item(132, 14).Value = -9823.842199999999  # N132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(132, 8).Value = 63678244  # H132
$ws.Cells.Item(132, 9).Value = 83705016  # I132
$ws.Cells.Item(132, 11).Value = 251115048  # K132
$ws.Cells.Item(132, 13).Value = -251112518  # M132

$ws.Cells.Item(136, 8).Value = 92403.09  # H136
$ws.Cells.Item(136, 9).Value = 112548.22  # I136
$ws.Cells.Item(136, 10).Value = 1750  # J136
$ws.Cells.Item(136, 11).Value = 337644.66  # K136
$ws.Cells.Item(136, 12).Value = 5250  # L136
$ws.Cells.Item(136, 13).Value = -335094.66  # M136
$ws.Cells.Item(136, 14).Value = -10350  # N136
